$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the price/volume columns so values like "1.00"
# or "224.78" are stored as literal text, matching the source data
# (these columns hold dotted/thousand-separated price strings, not numbers).
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "34.498.04"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.808.42"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "224.78"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  +5.38%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "38.80"
$ws.Range("E8").Value = "  +8.30%  "
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "2.067.39"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "1.802.38"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "34.445.98"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "68.18"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "241.50"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "11.10"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "171.07"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "7.68"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "17.59"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").Value = "0.0515"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Value = "1.317.04"
$ws.Range("E35").Value = "  -5.73%  "
$ws.Range("D36").Value = "0.640"
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "0.0186"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("E39").Value = "  -4.76%  "
$ws.Range("D40").Value = "82.46"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "0.947"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "13.74"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "1.968.39"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "5.79"
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "102.61"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "0.0₆0119"
$ws.Range("E51").Value = "  -6.34%  "

# Restore the default (unstyled) cell format so the sheet formatting
# matches the original, untouched cells.
$priceRange.Style = "Normal"
